$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3605.4773
$ws.Range("I15").Value = 3605.4773
$ws.Range("K15").Value = 10816.4319
$ws.Range("M15").Value = -10647.4319
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H113").Value = 3556.7778
$ws.Range("I113").Value = 2005
$ws.Range("J113").Value = 3750.75
$ws.Range("K113").Value = 2005
$ws.Range("L113").Value = 3750.75
$ws.Range("M113").Value = 1249
$ws.Range("N113").Value = -10258.75
$ws.Range("H132").Value = 9015949
$ws.Range("I132").Value = 9809251
$ws.Range("K132").Value = 29427753
$ws.Range("M132").Value = -29425223
$ws.Range("H137").Value = 1058.6296
$ws.Range("I137").Value = 950.125
$ws.Range("J137").Value = 1926.6666
$ws.Range("K137").Value = 2850.375
$ws.Range("L137").Value = 5779.9998
$ws.Range("M137").Value = -300.375
$ws.Range("N137").Value = -10879.9998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 216.66667
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 150
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = -38
$ws.Range("N5").Value = -474
$ws.Range("H8").Value = 3336001.8
$ws.Range("I8").Value = 5001502.5
$ws.Range("J8").Value = 5000
$ws.Range("K8").Value = 5001502.5
$ws.Range("L8").Value = 5000
$ws.Range("M8").Value = -5001358.5
$ws.Range("N8").Value = -5288
$ws.Range("H122").Value = 1327.7142
$ws.Range("I122").Value = 1339.5
$ws.Range("K122").Value = 4018.5
$ws.Range("M122").Value = -1568.5
$ws.Range("H132").Value = 2061.111
$ws.Range("I132").Value = 1674.4546
$ws.Range("K132").Value = 5023.3638
$ws.Range("M132").Value = -2493.3638
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 216.66667
$ws.Range("I4").Value = 150
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = -35
$ws.Range("N4").Value = -480
$ws.Range("H11").Value = 1369.6
$ws.Range("I11").Value = 2749.5
$ws.Range("J11").Value = 449.66666
$ws.Range("K11").Value = 2749.5
$ws.Range("L11").Value = 449.66666
$ws.Range("M11").Value = -2609.5
$ws.Range("N11").Value = -729.66666
$ws.Range("H86").Value = 3453.3928
$ws.Range("I86").Value = 3759.45
$ws.Range("J86").Value = 2688.25
$ws.Range("K86").Value = 3759.45
$ws.Range("L86").Value = 2688.25
$ws.Range("M86").Value = -2636.45
$ws.Range("N86").Value = -4934.25
$ws.Range("H89").Value = 3453.3928
$ws.Range("I89").Value = 3759.45
$ws.Range("J89").Value = 2688.25
$ws.Range("K89").Value = 18797.25
$ws.Range("L89").Value = 13441.25
$ws.Range("M89").Value = -13181.25
$ws.Range("N89").Value = -24673.25
$ws.Range("H104").Value = 55228
$ws.Range("J104").Value = 55228
$ws.Range("L104").Value = 55228
$ws.Range("N104").Value = -62216
$ws.Range("H134").Value = 8431.263000000001
$ws.Range("I134").Value = 1206.625
$ws.Range("K134").Value = 3619.875
$ws.Range("M134").Value = -1084.875
$ws.Range("H135").Value = 53267.5
$ws.Range("J135").Value = 53267.5
$ws.Range("L135").Value = 53267.5
$ws.Range("N135").Value = -63407.5
$ws.Range("H140").Value = 27488
$ws.Range("J140").Value = 27488
$ws.Range("L140").Value = 27488
$ws.Range("N140").Value = -37848
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 106.85714
$ws.Range("I7").Value = 67
$ws.Range("J7").Value = 206.5
$ws.Range("K7").Value = 67
$ws.Range("L7").Value = 206.5
$ws.Range("M7").Value = 46
$ws.Range("N7").Value = -432.5
$ws.Range("H31").Value = 1349.7931
$ws.Range("I31").Value = 966
$ws.Range("J31").Value = 2821
$ws.Range("K31").Value = 966
$ws.Range("L31").Value = 2821
$ws.Range("M31").Value = -671
$ws.Range("N31").Value = -3411
$ws.Range("H34").Value = 1349.7931
$ws.Range("I34").Value = 966
$ws.Range("J34").Value = 2821
$ws.Range("K34").Value = 966
$ws.Range("L34").Value = 2821
$ws.Range("M34").Value = -764
$ws.Range("N34").Value = -3225
$ws.Range("H122").Value = 714.7586
$ws.Range("I122").Value = 640.7917
$ws.Range("J122").Value = 1069.8
$ws.Range("K122").Value = 1922.3751
$ws.Range("L122").Value = 3209.4
$ws.Range("M122").Value = 527.6249
$ws.Range("N122").Value = -8109.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1684.9231
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 1684.9231
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 5054.7693
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -5642.7693
$ws.Range("H103").Value = 8378.857
$ws.Range("J103").Value = 14274.875
$ws.Range("L103").Value = 42824.625
$ws.Range("N103").Value = -44582.625
$ws.Range("H117").Value = 416.66666
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H121").Value = 551.1875
$ws.Range("I121").Value = 424.875
$ws.Range("J121").Value = 677.5
$ws.Range("K121").Value = 1274.625
$ws.Range("L121").Value = 2032.5
$ws.Range("M121").Value = 35.375
$ws.Range("N121").Value = -4652.5
$ws.Range("H129").Value = 41667096
$ws.Range("I129").Value = 37037510
$ws.Range("J129").Value = 83333336
$ws.Range("K129").Value = 111112530
$ws.Range("L129").Value = 250000008
$ws.Range("M129").Value = -111107530
$ws.Range("N129").Value = -250010008
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1371.96
$ws.Range("I102").Value = 1415.8889
$ws.Range("K102").Value = 1415.8889
$ws.Range("M102").Value = 206.1111000000001
$ws.Range("H107").Value = 677.52
$ws.Range("I107").Value = 778.9286
$ws.Range("J107").Value = 548.4545000000001
$ws.Range("K107").Value = 778.9286
$ws.Range("L107").Value = 548.4545000000001
$ws.Range("M107").Value = 1141.0714
$ws.Range("N107").Value = -4388.4545
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 300
$ws.Range("I9").Value = 300
$ws.Range("K9").Value = 300
$ws.Range("M9").Value = -160
$ws.Range("H81").Value = 325
$ws.Range("I81").Value = 300
$ws.Range("J81").Value = 350
$ws.Range("K81").Value = 600
$ws.Range("L81").Value = 700
$ws.Range("M81").Value = 461
$ws.Range("N81").Value = -2822
$ws.Range("H84").Value = 325
$ws.Range("I84").Value = 300
$ws.Range("J84").Value = 350
$ws.Range("K84").Value = 3000
$ws.Range("L84").Value = 3500
$ws.Range("M84").Value = 2304
$ws.Range("N84").Value = -14108
$ws.Range("H107").Value = 380
$ws.Range("I107").Value = 514
$ws.Range("J107").Value = 219.2
$ws.Range("K107").Value = 1542
$ws.Range("L107").Value = 657.5999999999999
$ws.Range("M107").Value = 378
$ws.Range("N107").Value = -4497.6
$ws.Range("H113").Value = 486.08334
$ws.Range("I113").Value = 284.6
$ws.Range("J113").Value = 1493.5
$ws.Range("K113").Value = 853.8000000000001
$ws.Range("L113").Value = 4480.5
$ws.Range("M113").Value = 1316.2
$ws.Range("N113").Value = -8820.5
